# Gallery_Sounder_FIM.xlsx - add "Portugal" market test data sheet
# (mirrors the "Swiss" sheet, with Portugal-specific values)

$wb = $excel.ActiveWorkbook

# --- Create the new "Portugal" sheet by copying "Swiss" (last tab) ---
$swiss = $wb.Worksheets.Item("Swiss")
$swiss.Copy($null, $swiss) | Out-Null
$portugal = $wb.Worksheets.Item($wb.Worksheets.Count)
$portugal.Name = "Portugal"

# --- Update the market-specific values ---
$portugal.Range("B2").Value = "Portugal Market"
$portugal.Range("B4").Value = "NGC-3479/T2438/T2465"

# --- Match column widths / row heights used on the Portugal sheet ---
$portugal.Columns.Item(1).ColumnWidth = 21.333333333333332
$portugal.Columns.Item(2).ColumnWidth = 19.666666666666668
$portugal.Columns.Item(3).ColumnWidth = 12
$portugal.Columns.Item(4).ColumnWidth = 11.833333333333334

$portugal.Rows.Item(3).RowHeight = 28.8
$portugal.Rows.Item(4).RowHeight = 28.8
$portugal.Rows.Item(5).RowHeight = 28.8

# --- Germany sheet: selection becomes the whole A1:D12 block ---
$germany = $wb.Worksheets.Item("Germany")
$germany.Activate()
$germany.Range("A1:D12").Select() | Out-Null

# --- Portugal becomes the active/selected tab, with B4 selected ---
$portugal.Activate()
$portugal.Range("B4").Select() | Out-Null
